$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cross table")

# ---------------------------------------------------------------------------
# Type in the cross-table text labels first (row by row, the way a person
# building this table in Excel would), so the shared-string table ends up
# in the same order as the source file.
# ---------------------------------------------------------------------------
$ws.Range("B15").Value = "Age Cohort"
$ws.Range("B16").Value = "18 to 25"
$ws.Range("C15").Value = "Employed"
$ws.Range("D15").Value = "Unemployed"
$ws.Range("B17").Value = "25 to 35"
$ws.Range("B18").Value = "35 to 45"
$ws.Range("B19").Value = "45 to 55"
$ws.Range("B20").Value = "55 to 65"
$ws.Range("B21").Value = "65+"
$ws.Range("B22").Value = "Total"
$ws.Range("E15").Value = "Total"

# ---------------------------------------------------------------------------
# Data rows (16-21): % employed, % unemployed, total (each row has one
# literal percentage coming from the background data, the other two cells
# are derived formulas)
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 0.6
$ws.Range("D16").Formula = "=100%-C16"
$ws.Range("E16").Formula = "=SUM(C16:D16)"

$ws.Range("C17").Value = 0.85
$ws.Range("D17").Formula = "=100%-C17"
$ws.Range("E17").Formula = "=SUM(C17:D17)"

$ws.Range("C18").Formula = "=100%-D18"
$ws.Range("D18").Value = 0.05
$ws.Range("E18").Formula = "=SUM(C18:D18)"

$ws.Range("C19").Formula = "=100%-D19"
$ws.Range("D19").Value = 0.03
$ws.Range("E19").Formula = "=SUM(C19:D19)"

$ws.Range("C20").Formula = "=100%-D20"
$ws.Range("D20").Value = 0.03
$ws.Range("E20").Formula = "=SUM(C20:D20)"

$ws.Range("C21").Value = 1
$ws.Range("D21").Formula = "=100%-C21"
$ws.Range("E21").Formula = "=SUM(C21:D21)"

# ---------------------------------------------------------------------------
# Total row (22)
# ---------------------------------------------------------------------------
$ws.Range("C22").Formula = "=SUM(C16:C21)"
$ws.Range("D22").Formula = "=SUM(D16:D21)"
$ws.Range("E22").Formula = "=SUM(E16:E21)"

# ---------------------------------------------------------------------------
# Number formatting: percentages for the data + total rows
# ---------------------------------------------------------------------------
$ws.Range("C16:E22").Style = "Percent"
$ws.Range("C16:E22").NumberFormat = "0%"

# ---------------------------------------------------------------------------
# Header + total row bold styling (matches the "Background"/"Task" labels'
# bold navy font already used elsewhere in the workbook)
# ---------------------------------------------------------------------------
$hdr = $ws.Range("B15:E15,B22:E22").Font
$hdr.Bold = $true
$hdr.Size = 10
$hdr.Color = 6299648
$hdr.Name = "Arial"

# Row labels (18 to 25 ... 65+) reuse the same bold navy style used for the
# other row labels already present in the sheet.
$lbl = $ws.Range("B16:B21").Font
$lbl.Bold = $true
$lbl.Size = 9
$lbl.Color = 6299648
$lbl.Name = "Arial"

# ---------------------------------------------------------------------------
# Column widths to fit the new headers / values
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 13.451822916666666
$ws.Columns.Item(4).ColumnWidth = 11.022135416666666
$ws.Columns.Item(5).ColumnWidth = 8.166666666666666

Write-Host "cross table created"
